# Apply the "B suite" update: append TestCase_B93..B98 rows (94-99) to the
# "Test Cases" sheet, fix D93's border style, resize the two wrapped rows,
# and leave the selection on the newly added content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Fix D93: border-only style (s="3") instead of fill+border (s="7") ---
$ws.Range("D92").Copy()
$ws.Range("D93").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 94 : TestCase_B93 ---
$ws.Range("A94").Value2 = "TestCase_B93"
$ws.Range("B94").Value2 = "OPQA-1238"
$ws.Range("C94").Value2 = "Verify that profile page of a person gets displayed when clicks on any PEOPLE search result in ALL search results page"
$ws.Range("D94").Value2 = "Y"
$ws.Range("E94").Value2 = "PASS"

# --- Row 95 : TestCase_B94 ---
$ws.Range("A95").Value2 = "TestCase_B94"
$ws.Range("B95").Value2 = "OPQA-1239"
$ws.Range("C95").Value2 = "Verify that record view page of a person gets displayed when user clicks on any PEOPLE in PEOPLE search results page."
$ws.Range("D95").Value2 = "Y"
$ws.Range("E95").Value2 = "FAIL"

# --- Row 96 : TestCase_B95 (tall wrapped row, result left blank) ---
$ws.Range("A96").Value2 = "TestCase_B95"
$ws.Range("B96").Value2 = "OPQA-599"
$ws.Range("C96").Value2 = "Verify that following fields get displayed correctly for a post in ALL search results page: `na)Title `nb)Creation date and time `nc)Author `nd)Author details `ne)Likes count `nf)Comments count"
$ws.Range("D96").Value2 = "Y"

# --- Row 97 : TestCase_B96 (tall wrapped row, result left blank) ---
$ws.Range("A97").Value2 = "TestCase_B96"
$ws.Range("B97").Value2 = "OPQA-553"
$ws.Range("C97").Value2 = "Verify that following fields get displayed correctly for a post in POSTS search results page: `na)Title `nb)Creation date and time `nc)Author `nd)Author details `ne)Likes count `nf)Comments count"
$ws.Range("D97").Value2 = "Y"

# --- Row 98 : TestCase_B97 ---
$ws.Range("A98").Value2 = "TestCase_B97"
$ws.Range("B98").Value2 = "OPQA-565"
$ws.Range("C98").Value2 = "Verify that no filtering options are present in ALL search results page"
$ws.Range("D98").Value2 = "Y"
$ws.Range("E98").Value2 = "PASS"

# --- Row 99 : TestCase_B98 ---
$ws.Range("A99").Value2 = "TestCase_B98"
$ws.Range("B99").Value2 = "OPQA-571"
$ws.Range("C99").Value2 = "Verify that search drop down content type is retained when user navigates back to ALL search results page from record view page"
$ws.Range("D99").Value2 = "Y"
$ws.Range("E99").Value2 = "PASS"

# --- Formatting: match the borders/wrap/fill already used on neighbouring rows ---
$ws.Range("A8").Copy()
$ws.Range("A94:A99").PasteSpecial(-4122)

$ws.Range("B71").Copy()
$ws.Range("B94:B99").PasteSpecial(-4122)

# C94/C95 are plain bordered cells (no wrap) like C26; C96-C99 wrap like C27
$ws.Range("C26").Copy()
$ws.Range("C94:C95").PasteSpecial(-4122)

$ws.Range("C27").Copy()
$ws.Range("C96:C99").PasteSpecial(-4122)

# D94 & D99 use the plain border style (s=3), the rest use fill+border (s=7)
$ws.Range("D2").Copy()
$ws.Range("D94").PasteSpecial(-4122)
$ws.Range("D99").PasteSpecial(-4122)

$ws.Range("D93").Copy()
$ws.Range("D95:D98").PasteSpecial(-4122)

# E94/E95/E98 use fill+border (s=7); E96/E97 plain+blank (s=3); E99 plain (s=3)
$ws.Range("E93").Copy()
$ws.Range("E94").PasteSpecial(-4122)
$ws.Range("E95").PasteSpecial(-4122)
$ws.Range("E98").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("E96").PasteSpecial(-4122)
$ws.Range("E97").PasteSpecial(-4122)
$ws.Range("E99").PasteSpecial(-4122)

# --- Row heights for the two wrapped description rows ---
$ws.Rows.Item(96).RowHeight = 105
$ws.Rows.Item(97).RowHeight = 105

# --- Move the view/selection to the freshly-entered data ---
[void]$ws.Range("A88").Select()
[void]$ws.Range("C96").Select()
